$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (cluster 0)
$ws.Range("B2").Value = 19950
$ws.Range("C2").Value = 79985
$ws.Range("D2").Value = 2019
$ws.Range("E2").Value = 112
$ws.Range("F2").Value = 159

# Row 3 (cluster 1)
$ws.Range("B3").Value = 11775
$ws.Range("C3").Value = 82698
$ws.Range("D3").Value = 2017
$ws.Range("E3").Value = 77
$ws.Range("F3").Value = 1010

# Row 4 (cluster 2)
$ws.Range("B4").Value = 23885
$ws.Range("C4").Value = 68783.5
$ws.Range("D4").Value = 2018
$ws.Range("E4").Value = 125
$ws.Range("F4").Value = 1782

# Row 5 (cluster 3)
$ws.Range("B5").Value = 32990
$ws.Range("C5").Value = 58000
$ws.Range("D5").Value = 2017
$ws.Range("E5").Value = 145
$ws.Range("F5").Value = 133

# Row 6 (cluster 4)
$ws.Range("B6").Value = 13990
$ws.Range("C6").Value = 89000
$ws.Range("D6").Value = 2016.5
$ws.Range("E6").Value = 110
$ws.Range("F6").Value = 150
